$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "Prioce" -> "Price" in the Price Grab description
$ws.Range("B2").Value = "Price Grab - Competetive Price scraping from Portals"

# Rename "DataWare house " -> "Data Warehouse "
$ws.Range("B6").Value = "Data Warehouse "

# Update the view: scroll position (top-left visible cell) and selected cell
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B11").Select()

